$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 corresponds to the 53f52a85-...zh-cn handback entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 10:14:42"
$wsZhCn.Range("G5").Value = "2016-02-17 10:15:30"

# de-de sheet: row 5 corresponds to the 53f52a85-...de-de handback entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 10:14:54"
$wsDeDe.Range("G5").Value = "2016-02-17 10:15:52"
